$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 214285.7142857143
$ws.Range("B13").Value = 142857.1428571429
$ws.Range("B23").Value = 214285.7142857143
$ws.Range("B31").Value = 214285.7142857143
$ws.Range("B32").Value = 142857.1428571429
$ws.Range("B33").Value = 214285.7142857143
$ws.Range("A34").Value = "Tổng lương tại HỆ THỐNG"
$ws.Range("B34").Value = 571428.5714285715
